$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a fresh "schema summary" row (used by the sheets that gain a
# brand-new row 3 describing a single named schema: 429, 500, 204, 401, 403,
# 404).
# Columns: A Section, B Name, C Parent, D Description, E Type, F Items Data
# Type, G Schema Name, H Format, I Mandatory, J Min, K Max, L PatternEba,
# M Regex, N Allowed value, O Example
# (Positional params only -- named params silently drop COM object values
# in this interpreter.)
# ---------------------------------------------------------------------------
function Set-SchemaSummaryRow3($ws, $section, $name) {
    $ws.Range("A3").Value = $section
    $ws.Range("B3").Value = $name
    $ws.Range("C3").Value = ""
    $ws.Range("D3").Value = ""
    $ws.Range("E3").Value = "schema"
    $ws.Range("F3").Value = ""
    $ws.Range("G3").Value = $name
    $ws.Range("H3").Value = ""
    $ws.Range("I3").Value = "Yes"
    $ws.Range("J3").Value = ""
    $ws.Range("K3").Value = ""
    $ws.Range("L3").Value = ""
    $ws.Range("M3").Value = ""
    $ws.Range("N3").Value = ""
    $ws.Range("O3").Value = ""
}

# ---------------------------------------------------------------------------
# Helper: collapse an existing multi-row body/content definition down to a
# single schema-reference row 3 (used by Body, 200 and 400), then drop the
# now-superfluous trailing rows.
# ---------------------------------------------------------------------------
function Set-SchemaReferenceRow3($ws, $name, $lastRow) {
    # A3 (Section) and C3/F3/H3/J3/K3/M3/N3 stay as they were.
    $ws.Range("B3").Value = $name
    $ws.Range("D3").Value = ""
    $ws.Range("E3").Value = "schema"
    $ws.Range("G3").Value = $name
    $ws.Range("I3").Value = "Yes"
    $ws.Range("L3").Value = ""
    $ws.Range("O3").Value = ""

    if ($lastRow -gt 3) {
        $ws.Rows("4:$lastRow").Delete()
    }
}

# --- "429" sheet: add errorResponse1 schema row -----------------------------
$ws = $wb.Worksheets.Item("429")
Set-SchemaSummaryRow3 $ws "content" "errorResponse1"

# --- "500" sheet: add errorResponse1 schema row -----------------------------
$ws = $wb.Worksheets.Item("500")
Set-SchemaSummaryRow3 $ws "content" "errorResponse1"

# --- "Body" sheet: collapse to a single request-schema row ------------------
$ws = $wb.Worksheets.Item("Body")
Set-SchemaReferenceRow3 $ws "interestMonthlyReport.211207Request" 7

# --- "200" sheet: collapse to a single response-schema row ------------------
$ws = $wb.Worksheets.Item("200")
Set-SchemaReferenceRow3 $ws "interestMonthlyReport.211207Response" 7

# --- "204" sheet: add interestMonthlyReport.211207Response schema row -------
$ws = $wb.Worksheets.Item("204")
Set-SchemaSummaryRow3 $ws "content" "interestMonthlyReport.211207Response"

# --- "400" sheet: collapse to a single errorResponse-schema row -------------
$ws = $wb.Worksheets.Item("400")
Set-SchemaReferenceRow3 $ws "errorResponse" 6

# --- "401" sheet: add errorResponse1 schema row -----------------------------
$ws = $wb.Worksheets.Item("401")
Set-SchemaSummaryRow3 $ws "content" "errorResponse1"

# --- "403" sheet: add errorResponse1 schema row -----------------------------
$ws = $wb.Worksheets.Item("403")
Set-SchemaSummaryRow3 $ws "content" "errorResponse1"

# --- "404" sheet: add errorResponse1 schema row -----------------------------
$ws = $wb.Worksheets.Item("404")
Set-SchemaSummaryRow3 $ws "content" "errorResponse1"
